$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, $ref, $val)
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws1 = $wb.Sheets.Item("Atlas Lions")

# Text/time values that must remain stored as text (shared strings)
Set-TextValue $ws1 "C16" "47.018"
$ws1.Range("C20").Value = "01:04.030"
Set-TextValue $ws1 "K3" "26.291"
Set-TextValue $ws1 "K16" "47.321"
$ws1.Range("K17").Value = "01:01.548"

# Numeric cell updates
$ws1.Range("B2").Value = 479
$ws1.Range("D2").Value = 1281
$ws1.Range("F2").Value = 3037
$ws1.Range("H2").Value = 1112
$ws1.Range("J2").Value = 3037
$ws1.Range("L2").Value = 957
$ws1.Range("D3").Value = 1346
$ws1.Range("F3").Value = 5792
$ws1.Range("H3").Value = 1261
$ws1.Range("J3").Value = 2080
$ws1.Range("L3").Value = 924
$ws1.Range("B4").Value = 1390
$ws1.Range("D4").Value = 2110
$ws1.Range("F4").Value = 5488
$ws1.Range("H4").Value = 649
$ws1.Range("J4").Value = 1496
$ws1.Range("L4").Value = 1178
$ws1.Range("B5").Value = 2198
$ws1.Range("D5").Value = 4783
$ws1.Range("F5").Value = 4485
$ws1.Range("H5").Value = 1504
$ws1.Range("J5").Value = 1431
$ws1.Range("L5").Value = 1711
$ws1.Range("B6").Value = 371
$ws1.Range("D6").Value = 2438
$ws1.Range("F6").Value = 5191
$ws1.Range("H6").Value = 749
$ws1.Range("J6").Value = 2550
$ws1.Range("L6").Value = 1186
$ws1.Range("B7").Value = 663
$ws1.Range("D7").Value = 2345
$ws1.Range("F7").Value = 3866
$ws1.Range("H7").Value = 2200
$ws1.Range("J7").Value = 1421
$ws1.Range("L7").Value = 1428
$ws1.Range("B8").Value = 1609
$ws1.Range("F8").Value = 3830
$ws1.Range("H8").Value = 1050
$ws1.Range("J8").Value = 1309
$ws1.Range("L8").Value = 1322
$ws1.Range("D9").Value = 2485
$ws1.Range("H9").Value = 1481
$ws1.Range("J9").Value = 1916
$ws1.Range("L9").Value = 1197
$ws1.Range("B10").Value = 208
$ws1.Range("D10").Value = 1599
$ws1.Range("F10").Value = 5854
$ws1.Range("H10").Value = 1091
$ws1.Range("J10").Value = 1140
$ws1.Range("L10").Value = 813
$ws1.Range("B11").Value = 155
$ws1.Range("D11").Value = 2476
$ws1.Range("F11").Value = 5813
$ws1.Range("H11").Value = 1994
$ws1.Range("J11").Value = 1181
$ws1.Range("L11").Value = 1110
$ws1.Range("D12").Value = 1736
$ws1.Range("H12").Value = 1376
$ws1.Range("J12").Value = 1444
$ws1.Range("L12").Value = 985
$ws1.Range("B13").Value = 502
$ws1.Range("D13").Value = 2648
$ws1.Range("H13").Value = 1748
$ws1.Range("J13").Value = 760
$ws1.Range("L13").Value = 1003
$ws1.Range("B14").Value = 305
$ws1.Range("D14").Value = 1386
$ws1.Range("F14").Value = 8327
$ws1.Range("H14").Value = 1458
$ws1.Range("J14").Value = 933
$ws1.Range("L14").Value = 874
$ws1.Range("B15").Value = 145
$ws1.Range("D15").Value = 1503
$ws1.Range("F15").Value = 2118
$ws1.Range("H15").Value = 1934
$ws1.Range("J15").Value = 888
$ws1.Range("L15").Value = 845
$ws1.Range("B16").Value = 953
$ws1.Range("D16").Value = 924
$ws1.Range("F16").Value = 1393
$ws1.Range("H16").Value = 851
$ws1.Range("J16").Value = 2150
$ws1.Range("L16").Value = 909
$ws1.Range("B17").Value = 531
$ws1.Range("D17").Value = 6611
$ws1.Range("F17").Value = 8638
$ws1.Range("H17").Value = 796
$ws1.Range("J17").Value = 1349
$ws1.Range("L17").Value = 892
$ws1.Range("B18").Value = 610
$ws1.Range("D18").Value = 7418
$ws1.Range("F18").Value = 3157
$ws1.Range("H18").Value = 839
$ws1.Range("J18").Value = 1760
$ws1.Range("L18").Value = 1069
$ws1.Range("B19").Value = 227
$ws1.Range("D19").Value = 2599
$ws1.Range("F19").Value = 6036
$ws1.Range("H19").Value = 3059
$ws1.Range("J19").Value = 1190
$ws1.Range("L19").Value = 1338
$ws1.Range("B20").Value = 317
$ws1.Range("D20").Value = 1758
$ws1.Range("F20").Value = 5014
$ws1.Range("H20").Value = 1081
$ws1.Range("J20").Value = 965
$ws1.Range("L20").Value = 787
$ws1.Range("D21").Value = 3304
$ws1.Range("F21").Value = 2941
$ws1.Range("H21").Value = 1037
$ws1.Range("L21").Value = 566

$ws2 = $wb.Sheets.Item("Austria's Bench Team YEP")

# Text/time values that must remain stored as text (shared strings)
Set-TextValue $ws2 "C2" "35.115"
Set-TextValue $ws2 "C3" "25.943"
Set-TextValue $ws2 "E3" "26.355"
Set-TextValue $ws2 "G8" "39.383"
$ws2.Range("K20").Value = "01:11.988"

# Numeric cell updates
$ws2.Range("B2").Value = 1239
$ws2.Range("D2").Value = 5032
$ws2.Range("F2").Value = 168
$ws2.Range("H2").Value = 2634
$ws2.Range("J2").Value = 5209
$ws2.Range("L2").Value = 1347
$ws2.Range("B3").Value = 526
$ws2.Range("D3").Value = 2621
$ws2.Range("F3").Value = 834
$ws2.Range("H3").Value = 546
$ws2.Range("L3").Value = 635
$ws2.Range("B4").Value = 782
$ws2.Range("D4").Value = 2664
$ws2.Range("H4").Value = 1555
$ws2.Range("J4").Value = 1743
$ws2.Range("L4").Value = 823
$ws2.Range("B5").Value = 339
$ws2.Range("D5").Value = 3280
$ws2.Range("F5").Value = 141
$ws2.Range("H5").Value = 8394
$ws2.Range("L5").Value = 1253
$ws2.Range("B6").Value = 1926
$ws2.Range("D6").Value = 1991
$ws2.Range("F6").Value = 333
$ws2.Range("H6").Value = 306
$ws2.Range("L6").Value = 855
$ws2.Range("B7").Value = 2730
$ws2.Range("D7").Value = 2673
$ws2.Range("F7").Value = 912
$ws2.Range("H7").Value = 1459
$ws2.Range("L7").Value = 1681
$ws2.Range("D8").Value = 2846
$ws2.Range("F8").Value = 70
$ws2.Range("H8").Value = 193
$ws2.Range("J8").Value = 4632
$ws2.Range("L8").Value = 551
$ws2.Range("B9").Value = 1607
$ws2.Range("D9").Value = 1935
$ws2.Range("F9").Value = 641
$ws2.Range("H9").Value = 468
$ws2.Range("J9").Value = 2965
$ws2.Range("L9").Value = 905
$ws2.Range("B10").Value = 3196
$ws2.Range("D10").Value = 2626
$ws2.Range("F10").Value = 706
$ws2.Range("H10").Value = 729
$ws2.Range("J10").Value = 2814
$ws2.Range("L10").Value = 1353
$ws2.Range("B11").Value = 3984
$ws2.Range("D11").Value = 4123
$ws2.Range("F11").Value = 655
$ws2.Range("H11").Value = 1164
$ws2.Range("J11").Value = 6053
$ws2.Range("L11").Value = 1934
$ws2.Range("B12").Value = 3415
$ws2.Range("D12").Value = 1684
$ws2.Range("H12").Value = 1733
$ws2.Range("J12").Value = 7663
$ws2.Range("L12").Value = 1183
$ws2.Range("B13").Value = 2726
$ws2.Range("D13").Value = 1289
$ws2.Range("F13").Value = 923
$ws2.Range("J13").Value = 4084
$ws2.Range("L13").Value = 776
$ws2.Range("D14").Value = 1800
$ws2.Range("F14").Value = 670
$ws2.Range("H14").Value = 1891
$ws2.Range("L14").Value = 1453
$ws2.Range("B15").Value = 2347
$ws2.Range("D15").Value = 2283
$ws2.Range("F15").Value = 462
$ws2.Range("H15").Value = 6525
$ws2.Range("L15").Value = 1697
$ws2.Range("B16").Value = 9832
$ws2.Range("D16").Value = 2478
$ws2.Range("F16").Value = 851
$ws2.Range("H16").Value = 2750
$ws2.Range("L16").Value = 2026
$ws2.Range("B17").Value = 9159
$ws2.Range("D17").Value = 2671
$ws2.Range("F17").Value = 1002
$ws2.Range("H17").Value = 211
$ws2.Range("L17").Value = 1294
$ws2.Range("B18").Value = 5683
$ws2.Range("D18").Value = 3142
$ws2.Range("H18").Value = 2559
$ws2.Range("L18").Value = 1922
$ws2.Range("B19").Value = 7857
$ws2.Range("D19").Value = 2444
$ws2.Range("F19").Value = 920
$ws2.Range("H19").Value = 651
$ws2.Range("L19").Value = 1338
$ws2.Range("B20").Value = 7119
$ws2.Range("D20").Value = 1451
$ws2.Range("F20").Value = 270
$ws2.Range("H20").Value = 202
$ws2.Range("J20").Value = 19999
$ws2.Range("L20").Value = 641
$ws2.Range("B21").Value = 1948
$ws2.Range("D21").Value = 2245
$ws2.Range("F21").Value = 829
$ws2.Range("H21").Value = 819
$ws2.Range("L21").Value = 1198
